$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting existing columns B..W to C..X.
$ws.Columns("B:B").Insert()

# Populate the new header cell with the new option text.
$ws.Range("B1").Value = "申請年份 Year of Application"

# Match the saved selection/view state: B1 selected, no special scroll position.
$ws.Range("B1").Select()
